$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$fixedUrl = "https://github.com/nguyentienminh07102004/product-management/commit/f1ac942243f94ace9433ec5239d0ec416202bf6b"
$newUrl   = "https://github.com/nguyentienminh07102004/product-management/commit/a4a024e0a5282e29b9c0298532c0093a1674ed60"

# --- C12 formatting -------------------------------------------------------
# Give C12 the "highlighted hyperlink" look (fill + Hyperlink font) that the
# other link cells in the highlighted block use, mirrored from C10 which is
# already the exact combination required (fontId=Hyperlink, fillId=2,
# borderId=0, linked to the Hyperlink cell style).
$ws.Range("C10").Copy() | Out-Null
$ws.Range("C12").PasteSpecial(-4122) | Out-Null  # xlPasteFormats
$excel.CutCopyMode = 0

# --- Cell values ------------------------------------------------------------
$ws.Range("C5").Value = $fixedUrl
$ws.Range("C12").Value = $newUrl

# --- Hyperlinks --------------------------------------------------------------
# The engine only supports clearing the *entire* hyperlink collection at once,
# so capture every existing target, wipe the collection, then recreate all of
# them (fixing C5's target and adding the brand-new C12 one) in their original
# left-to-right, top-to-bottom order.
$links = @(
    @{ Cell = "C2";  Url = "https://github.com/nguyentienminh07102004/product-management/commit/0670d984c38014165261ac28c8d9bf6f6bda634d" },
    @{ Cell = "C3";  Url = "https://github.com/nguyentienminh07102004/product-management/commit/262d8a2429fc1cdc879187949eea553f2f23143c" },
    @{ Cell = "C5";  Url = $fixedUrl },
    @{ Cell = "C6";  Url = "https://github.com/nguyentienminh07102004/product-management/commit/41534278856704293c82600456b1f6467babb5b8" },
    @{ Cell = "C7";  Url = "https://github.com/nguyentienminh07102004/product-management/commit/23586a985d09cd41448981b336656f09c95dd7dd" },
    @{ Cell = "C8";  Url = "https://github.com/nguyentienminh07102004/product-management/commit/a36764f2bc798b74ffb3bdda7168d990212f567d" },
    @{ Cell = "C9";  Url = "https://github.com/nguyentienminh07102004/product-management/commit/e72573281984dbf2b78af3ac8215ef1773f841a9" },
    @{ Cell = "C4";  Url = "https://github.com/nguyentienminh07102004/product-management/commit/e896019a677ede65944dae61c70a40c4ecc8a67d" },
    @{ Cell = "C10"; Url = "https://github.com/nguyentienminh07102004/product-management/commit/9516b17e3bbe0f9823dc17f0bdf0c65b8008935b" },
    @{ Cell = "C11"; Url = "https://github.com/nguyentienminh07102004/product-management/commit/6cacbdd5ab0f00482423975c4dfea365d3835ef6" },
    @{ Cell = "C12"; Url = $newUrl }
)

$ws.Range("C2").Hyperlinks.Delete()

foreach ($l in $links) {
    $ws.Hyperlinks.Add($ws.Range($l.Cell), $l.Url) | Out-Null
}

# --- Selection ---------------------------------------------------------------
# Matches the final saved cursor position in the workbook.
$ws.Range("C7").Select()

$wb.Save()
